$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "45.899.93"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.579.36"
$ws.Range("E3").Value = "  +8.39%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.08"
$ws.Range("E5").Value = "  +1.05%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.99"
$ws.Range("E6").Value = "  +0.15%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +10.92%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.04"
$ws.Range("E10").Value = "  +10.09%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +5.25%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.02"
$ws.Range("E12").Value = "  +12.21%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.975.73"
$ws.Range("E13").Value = "  +8.50%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.86%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "2.583.85"
$ws.Range("E15").Value = "  +9.50%  "

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.888"
$ws.Range("E16").Value = "  +7.50%  "

# Row 17 - Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.65"
$ws.Range("E17").Value = "  +6.44%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "45.984.37"
$ws.Range("E18").Value = "  +0.03%  "

# Row 19 - ShibaInu (was InternetComputer(DFINITY))
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0999"
$ws.Range("E19").Value = "  +5.05%  "

# Row 20 - InternetComputer(DFINITY) (was ShibaInu)
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("E20").Value = "  +1.82%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +8.78%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  +5.39%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.76"
$ws.Range("E23").Value = "  +3.27%  "

# Row 24 - PancakeSwap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  +5.81%  "

# Row 25 - ImmutableX
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.42"
$ws.Range("E26").Value = "  +30.92%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.02%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +5.64%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.06"
$ws.Range("E29").Value = "  -2.23%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +1.93%  "

# Row 31 - Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  +8.80%  "

# Row 32 - LidoDAOToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.63"
$ws.Range("E32").Value = "  -3.01%  "

# Row 33 - WEMIXToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  +4.21%  "

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.23"
$ws.Range("E34").Value = "  +16.46%  "

# Row 35 - Monero
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.66"
$ws.Range("E35").Value = "  +3.19%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0822"
$ws.Range("E36").Value = "  +6.40%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +1.66%  "

# Row 38 - Stellar
$ws.Range("E38").Value = "  +4.51%  "

# Row 39 - Celestia
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.76"
$ws.Range("E39").Value = "  +4.59%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.11"
$ws.Range("E40").Value = "  +5.92%  "

# Row 41 - NEARProtocol
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  +9.43%  "

# Row 42 - VeChain
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0317"
$ws.Range("E42").Value = "  +5.97%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.043.31"
$ws.Range("E43").Value = "  +5.62%  "

# Row 44 - EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.57"
$ws.Range("E44").Value = "  +38.12%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46 - BitcoinSV
$ws.Range("E46").Value = "  -1.70%  "

# Row 47 - FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  +7.12%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -2.04%  "

# Row 49 - Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.99"
$ws.Range("E49").Value = "  +9.11%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.831.36"
$ws.Range("E50").Value = "  +8.35%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +6.18%  "
